# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
# Concretely:
#   1. Insert a new "Player Info" worksheet as the first sheet, with the
#      player's ID / NAME / BATTING_HAND / BOWL_STYLE.
#   2. Rename the MATCH_CARD_LINK column to MATCH_CODE on both the
#      "ODI Batting" and "ODI Bowling" sheets, and replace the full
#      howstat.com scorecard URL with just the bare match code ("4521").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. New "Player Info" sheet (inserted at the front of the workbook -
#    Worksheets.Add() with no placement args inserts before the first
#    existing sheet, same as it would in Excel).
# ---------------------------------------------------------------------
$playerInfo = $wb.Worksheets.Add()
$playerInfo.Name = "Player Info"

# Header row - bold, centered, thin-bordered, matching the look of the
# other sheets' header rows.
$header = $playerInfo.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Data row. Force the ID cell to stay text (it is a player id, not a
# number) so it round-trips as "6466" rather than being coerced to a
# numeric cell.
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "6466"
$playerInfo.Range("B2").Value = "Chamika Deemantha Gunasekara"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium Fast"

# ---------------------------------------------------------------------
# 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE, full URL -> bare code
# ---------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"
$batting.Range("D2").NumberFormat = "@"
$batting.Range("D2").Value = "4521"

# ---------------------------------------------------------------------
# 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE, full URL -> bare code
# ---------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"
$bowling.Range("B2").NumberFormat = "@"
$bowling.Range("B2").Value = "4521"
